$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Durations_psplib")

$ws.Range("C1").Value = 2382.41220664978

$ws.Range("A4").Value = 304053.2375
$ws.Range("B4").Value = 306298
$ws.Range("F4").Value = 102251.8
$ws.Range("G4").Value = 103074
